$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "mean" to "mean score" in A3
$ws.Range("A3").Value = "mean score"

# Insert 3 new rows before the existing row 5 (p-value), pushing it down to row 8
$ws.Rows("5:7").Insert()

# Row 5: group total played
$ws.Range("A5").Value = "group total played"
$row5 = @(3890, 4378, 3868, 3958, 3889, 3974, 3952, 4001, 4159, 4318, 4063, 4237, 4425)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $row5[$i]
}

# Row 6: group mean
$ws.Range("A6").Value = "group mean"
$row6 = @(4.1159, 4.1039, 4.1107, 4.0942, 4.1273, 4.079, 4.1015, 4.1135, 4.076, 4.1061, 4.1088, 4.0991, 4.1024)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, $i + 2).Value = $row6[$i]
}

# Row 7: group st. dev.
$ws.Range("A7").Value = "group st. dev."
$row7 = @(1.0605, 1.056, 1.0586, 1.0519, 1.0613, 1.0461, 1.0554, 1.0603, 1.0512, 1.0569, 1.0599, 1.0591, 1.0562)
for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, $i + 2).Value = $row7[$i]
}
